$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 28.063402
$ws.Range("N2").Value = 84.190206
$ws.Range("O2").Value = 0.2422582722789452
$ws.Range("P2").Value = 0.2422582722789452
$ws.Range("Q2").Value = 3233.503606591805
$ws.Range("R2").Value = 29101.53245932625
$ws.Range("S2").Value = 0.06752728265886991
$ws.Range("T2").Value = 0.0675272826588699
$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.270119931190762
$ws.Range("P3").Value = 0.2701199311907619
$ws.Range("Q3").Value = 3605.382650099787
$ws.Range("R3").Value = 32448.44385089808
$ws.Range("S3").Value = 0.07529346582770279
$ws.Range("T3").Value = 0.07529346582770274
$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("M4").Value = 25.48508733333334
$ws.Range("N4").Value = 76.455262
$ws.Range("O4").Value = 0.22000088322333
$ws.Range("P4").Value = 0.2200008832233299
$ws.Range("Q4").Value = 2936.426660126255
$ws.Range("R4").Value = 26427.8399411363
$ws.Range("S4").Value = 0.06132323857043367
$ws.Range("T4").Value = 0.06132323857043365
$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 31.00143166666667
$ws.Range("N5").Value = 93.004295
$ws.Range("O5").Value = 0.2676209133069629
$ws.Range("P5").Value = 0.2676209133069628
$ws.Range("Q5").Value = 3572.027407927095
$ws.Range("R5").Value = 32148.24667134386
$ws.Range("S5").Value = 0.07459688739749518
$ws.Range("T5").Value = 0.07459688739749513
$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 28.063402
$ws.Range("N6").Value = 84.190206
$ws.Range("O6").Value = 0.2422582722789452
$ws.Range("P6").Value = 0.2422582722789452
$ws.Range("Q6").Value = 5187.289795930404
$ws.Range("R6").Value = 46685.60816337364
$ws.Range("S6").Value = 0.1083294243337714
$ws.Range("T6").Value = 0.1083294243337714
$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.270119931190762
$ws.Range("P7").Value = 0.2701199311907619
$ws.Range("S7").Value = 0.1207881835022744
$ws.Range("T7").Value = 0.1207881835022744
$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("M8").Value = 25.48508733333334
$ws.Range("N8").Value = 76.455262
$ws.Range("O8").Value = 0.22000088322333
$ws.Range("P8").Value = 0.2200008832233299
$ws.Range("Q8").Value = 4710.709466820708
$ws.Range("R8").Value = 42396.38520138637
$ws.Range("S8").Value = 0.09837669858828557
$ws.Range("T8").Value = 0.09837669858828554
$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 31.00143166666667
$ws.Range("N9").Value = 93.004295
$ws.Range("O9").Value = 0.2676209133069629
$ws.Range("P9").Value = 0.2676209133069628
$ws.Range("Q9").Value = 5730.36049384653
$ws.Range("R9").Value = 51573.24444461877
$ws.Range("S9").Value = 0.1196707101288986
$ws.Range("T9").Value = 0.1196707101288985
$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 28.063402
$ws.Range("N10").Value = 84.190206
$ws.Range("O10").Value = 0.2422582722789452
$ws.Range("P10").Value = 0.2422582722789452
$ws.Range("Q10").Value = 1699.265838421247
$ws.Range("R10").Value = 15293.39254579122
$ws.Range("S10").Value = 0.03548683364685625
$ws.Range("T10").Value = 0.03548683364685624
$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.270119931190762
$ws.Range("P11").Value = 0.2701199311907619
$ws.Range("Q11").Value = 1894.695140980133
$ws.Range("R11").Value = 17052.2562688212
$ws.Range("S11").Value = 0.03956810627225762
$ws.Range("T11").Value = 0.03956810627225761
$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("M12").Value = 25.48508733333334
$ws.Range("N12").Value = 76.455262
$ws.Range("O12").Value = 0.22000088322333
$ws.Range("P12").Value = 0.2200008832233299
$ws.Range("Q12").Value = 1543.146418766883
$ws.Range("R12").Value = 13888.31776890194
$ws.Range("S12").Value = 0.03222649394658578
$ws.Range("T12").Value = 0.03222649394658577
$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 31.00143166666667
$ws.Range("N13").Value = 93.004295
$ws.Range("O13").Value = 0.2676209133069629
$ws.Range("P13").Value = 0.2676209133069628
$ws.Range("Q13").Value = 1877.166345452961
$ws.Range("R13").Value = 16894.49710907665
$ws.Range("S13").Value = 0.03920204144776822
$ws.Range("T13").Value = 0.0392020414477682
$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 28.063402
$ws.Range("N14").Value = 84.190206
$ws.Range("O14").Value = 0.2422582722789452
$ws.Range("P14").Value = 0.2422582722789452
$ws.Range("Q14").Value = 1480.333464000883
$ws.Range("R14").Value = 13323.00117600795
$ws.Range("S14").Value = 0.03091473163944762
$ws.Range("T14").Value = 0.03091473163944762
$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.270119931190762
$ws.Range("P15").Value = 0.2701199311907619
$ws.Range("Q15").Value = 1650.583774389667
$ws.Range("R15").Value = 14855.253969507
$ws.Range("S15").Value = 0.03447017558852716
$ws.Range("T15").Value = 0.03447017558852716
$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("M16").Value = 25.48508733333334
$ws.Range("N16").Value = 76.455262
$ws.Range("O16").Value = 0.22000088322333
$ws.Range("P16").Value = 0.2200008832233299
$ws.Range("Q16").Value = 1344.328375174128
$ws.Range("R16").Value = 12098.95537656715
$ws.Range("S16").Value = 0.02807445211802496
$ws.Range("T16").Value = 0.02807445211802495
$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 31.00143166666667
$ws.Range("N17").Value = 93.004295
$ws.Range("O17").Value = 0.2676209133069629
$ws.Range("P17").Value = 0.2676209133069628
$ws.Range("Q17").Value = 1635.313378189264
$ws.Range("R17").Value = 14717.82040370338
$ws.Range("S17").Value = 0.0341512743328009
$ws.Range("T17").Value = 0.03415127433280089
